# Bot-Varredura/dominios.xlsx — "incluindo DTB - IBGE"
# Replace the placeholder domain with the generic prompt text, then lay out
# the sheet so it has room for the DTB/IBGE domain list (rows 2-94),
# matching column widths and row heights picked up when the list was pasted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Swap the sample domain for the placeholder text.
$ws.Range("A2").Value = "coloque_o_dominio_aqui"

# Give A2 an explicit (but visually-default) style - this is what happens
# when the column/row formatting used for the rest of the list is picked
# up by the first data cell too.
$ws.Range("A2").Style = "Normal"

# 2) Column widths for DOMINIO / STATUS / DATA EXTRACAO.
$ws.Columns.Item(1).ColumnWidth = 31
$ws.Columns.Item(2).ColumnWidth = 8.333333333333332
$ws.Columns.Item(3).ColumnWidth = 18.666666666666668

# 3) Row heights: header keeps its height, data rows get the tighter 13.8pt
# height used throughout the pasted list.
$ws.Rows.Item(2).RowHeight = 13.8
for ($r = 3; $r -le 93; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# 4) Touch the bottom-right corner cell so the sheet's used range grows to
# C94 (the full extent of the pasted DTB/IBGE domain table) without
# actually putting data/formatting into the still-empty rows.
$ws.Range("C94").Font.Bold = $false
$ws.Rows.Item(94).RowHeight = 13.8

# 5) Park the selection where the user left it after pasting the list.
$ws.Range("A3").Select()
